$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two more work-log entries were logged for 2014-07-26/27, so a new blank
# separator row is needed between the data block and the summary block.
# Insert it right above the old "sum [min]" row (159); this pushes the
# three summary rows down to 160/161/162 and keeps their formulas'
# relative cell references intact.
$ws.Rows.Item(159).Insert()

# --- Row 157: new data row (2014-07-26, 20:15-22:15) ---
$ws.Range("A157").Value = 2014
$ws.Range("B157").Value = 7
$ws.Range("C157").Value = 26
$ws.Range("D157").Value = 0.84375
$ws.Range("E157").Value = 0.92708333333333337

# --- Row 158: new data row (2014-07-27, 08:00-12:00) ---
$ws.Range("A158").Value = 2014
$ws.Range("B158").Value = 7
$ws.Range("C158").Value = 27
$ws.Range("D158").Value = 0.33333333333333331
$ws.Range("E158").Value = 0.5

# Fill the "time spent" formulas down into the two new rows, same pattern
# used for the rest of the table.
$ws.Range("F157:F158").Formula = "=(E157-D157)*24*60"
$ws.Range("G157:G158").Formula = "=F157/60"

# --- Summary block (now on rows 160-162) picks up the two extra rows ---
$ws.Range("F160").Formula = "=SUM(F2:F158)"
$ws.Range("F161").Formula = "=F160/60"
$ws.Range("F162").Formula = "=F161/38.5"

# --- Restore the active selection ---
$ws.Range("F158").Select()
